$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.360.17"
$ws.Range("E2").Value = "  -2.57%  "
$ws.Range("D3").Value = "1.938.26"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.43"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7248"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.46%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3335"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "28.39"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07241"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8105"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.79%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08099"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.95%  "
$ws.Range("D13").Value = "1.940.42"
$ws.Range("E13").Value = "  -2.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.473"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.54"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.36%  "
$ws.Range("D17").Value = "30.369.92"
$ws.Range("E17").Value = "  -2.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008238"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "249.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -8.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.923"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.96%  "
$ws.Range("D21").Value = "2.192.84"
$ws.Range("E21").Value = "  -2.53%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.939"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.761"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.97%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.51%  "
$ws.Range("E27").Value = "  +1.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1328"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.41%  "
$ws.Range("E30").Value = "  -1.72%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.346"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.79%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.444"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.200"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05198"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.291"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7519"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.751"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01980"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.839"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.517"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4542"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.040"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8484"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.65%  "
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.814"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.453"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.48%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4198"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06055"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.28%  "
